# Regenerate save_data to use K (strikeouts) instead of Strike# placeholder values.
# Column G ("K") values are recalculated/rewritten for rows 2-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 1
    4  = 2
    5  = 3
    6  = 2
    7  = 2
    8  = 0
    9  = 1
    10 = 0
    11 = 3
    12 = 0
    13 = 2
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 2
    19 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
